$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement of the crypto table rows 2-51 (B:E) per the day-over-day data refresh.
$data = @(
    ,@('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '27.113.95', '  -0.52%  ')
    ,@('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.894.83', '  -0.65%  ')
    ,@('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.003', '  +0.20%  ')
    ,@('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '306.79', '  -0.26%  ')
    ,@('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.002', '  +0.18%  ')
    ,@('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5216', '  -0.66%  ')
    ,@('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3762', '  -0.69%  ')
    ,@('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07269', '  +0.00%  ')
    ,@('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '21.16', '  -0.81%  ')
    ,@('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.8999', '  -0.08%  ')
    ,@('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.08209', '  +0.98%  ')
    ,@('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.963.12', '  +2.79%  ')
    ,@('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '96.20', '  +0.92%  ')
    ,@('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.302', '  +0.14%  ')
    ,@('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.003', '  +0.16%  ')
    ,@('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000008601', '  -0.27%  ')
    ,@('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '14.58', '  +0.57%  ')
    ,@('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.004', '  +0.30%  ')
    ,@('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '27.150.44', '  -0.63%  ')
    ,@('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.086', '  +0.40%  ')
    ,@('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.200.34', '  +2.23%  ')
    ,@('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.70', '  +0.68%  ')
    ,@('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.413', '  -0.74%  ')
    ,@('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.312', '  +0.06%  ')
    ,@('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '148.52', '  +1.50%  ')
    ,@('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.19', '  -0.04%  ')
    ,@('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.743', '  -0.24%  ')
    ,@('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '115.18', '  +0.11%  ')
    ,@('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.790', '  -0.50%  ')
    ,@('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.864', '  -2.56%  ')
    ,@('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.09202', '  -0.34%  ')
    ,@('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05024', '  -0.73%  ')
    ,@('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7923', '  -1.55%  ')
    ,@('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.215', '  -2.37%  ')
    ,@('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '3.429', '  +3.02%  ')
    ,@('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.964', '  -1.17%  ')
    ,@('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.613', '  +1.24%  ')
    ,@('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5722', '  -0.13%  ')
    ,@('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01994', '  +0.40%  ')
    ,@('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.074', '  -0.35%  ')
    ,@('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '9.025', '  +0.51%  ')
    ,@('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.550', '  -1.21%  ')
    ,@('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '116.47', '  -2.59%  ')
    ,@('Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1513', '  -0.13%  ')
    ,@('Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4855', '  +0.05%  ')
    ,@('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.002', '  +0.21%  ')
    ,@('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '10.11', '  -1.44%  ')
    ,@('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.620', '  -0.34%  ')
    ,@('Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '38.22', '  +1.55%  ')
    ,@('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '63.60', '  -0.32%  ')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}
